# DAN_YR_FIN.xlsx update — add the newest fiscal-year column (2018-12-31)
# as a new column D on the "DAN" sheet, pushing the existing D:K year
# columns one column to the right (to E:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank column at D; everything from D:K shifts to E:L.
$ws.Columns(4).Insert()

# New column D should carry the same number formatting as the data that
# just shifted into column E (date format for the header row, plain
# number format for the data rows), so copy E's formatting/values into D
# first and then overwrite with the real FY2018 figures below.
$ws.Range("E5:E102").Copy($ws.Range("D5:D102"))

$newColumnD = @(
    @{Row=7; Value=43465}
    @{Row=8; Value=8143000}
    @{Row=9; Value=6981000}
    @{Row=10; Value=1162000}
    @{Row=12; Value="NA"}
    @{Row=13; Value=0}
    @{Row=14; Value=63000}
    @{Row=15; Value=8000}
    @{Row=17; Value=7564000}
    @{Row=18; Value=579000}
    @{Row=20; Value=11000}
    @{Row=21; Value=860000}
    @{Row=22; Value=96000}
    @{Row=23; Value=494000}
    @{Row=24; Value=78000}
    @{Row=25; Value=0}
    @{Row=26; Value=416000}
    @{Row=27; Value=427000}
    @{Row=28; Value=0}
    @{Row=29; Value="NA"}
    @{Row=30; Value=0}
    @{Row=31; Value=0}
    @{Row=32; Value=-11000}
    @{Row=33; Value=427000}
    @{Row=34; Value=0}
    @{Row=35; Value=427000}
    @{Row=38; Value=43465}
    @{Row=41; Value=510000}
    @{Row=42; Value=21000}
    @{Row=43; Value=1243000}
    @{Row=44; Value=1031000}
    @{Row=45; Value=102000}
    @{Row=46; Value=2907000}
    @{Row=47; Value=233000}
    @{Row=48; Value=1850000}
    @{Row=49; Value=428000}
    @{Row=50; Value=0}
    @{Row=51; Value=0}
    @{Row=52; Value=500000}
    @{Row=53; Value=0}
    @{Row=54; Value=5918000}
    @{Row=57; Value=1217000}
    @{Row=58; Value=28000}
    @{Row=59; Value=502000}
    @{Row=60; Value=1747000}
    @{Row=61; Value=1755000}
    @{Row=62; Value=874000}
    @{Row=63; Value=0}
    @{Row=64; Value=0}
    @{Row=65; Value=0}
    @{Row=66; Value=4573000}
    @{Row=68; Value=0}
    @{Row=69; Value=0}
    @{Row=70; Value=0}
    @{Row=71; Value=0}
    @{Row=72; Value=456000}
    @{Row=73; Value=0}
    @{Row=74; Value=0}
    @{Row=75; Value=0}
    @{Row=76; Value=1345000}
    @{Row=77; Value=0}
    @{Row=80; Value=43465}
    @{Row=81; Value=427000}
    @{Row=83; Value=270000}
    @{Row=84; Value=0}
    @{Row=85; Value=0}
    @{Row=86; Value=0}
    @{Row=87; Value=0}
    @{Row=88; Value=0}
    @{Row=89; Value=568000}
    @{Row=91; Value=-325000}
    @{Row=92; Value=0}
    @{Row=93; Value=0}
    @{Row=94; Value=-462000}
    @{Row=96; Value=-58000}
    @{Row=97; Value=0}
    @{Row=98; Value=0}
    @{Row=99; Value=0}
    @{Row=100; Value=-180000}
    @{Row=101; Value=-16000}
    @{Row=102; Value=-90000}
)

foreach ($entry in $newColumnD) {
    $ws.Cells.Item($entry.Row, 4).Value = $entry.Value
}
